# Apply cell updates from the crypto price refresh.
# Values that look like plain numbers (e.g. "341.59") are prefixed with a
# leading apostrophe so Excel stores them as literal text, matching the
# original inline-string cell contents (e.g. "29.787.71", "1.037", ...).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.787.71'
$ws.Range('E2').Value = '  +7.00%  '
$ws.Range('D3').Value = '1.949.47'
$ws.Range('E3').Value = '  +5.36%  '
$ws.Range('E4').Value = '  -0.58%  '
$ws.Range('D5').Value = "'" + '341.59'
$ws.Range('E5').Value = '  +1.87%  '
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('D7').Value = "'" + '0.4793'
$ws.Range('E7').Value = '  +2.94%  '
$ws.Range('D8').Value = "'" + '0.4131'
$ws.Range('E8').Value = '  +6.80%  '
$ws.Range('D9').Value = "'" + '47.83'
$ws.Range('E9').Value = '  +2.29%  '
$ws.Range('D10').Value = "'" + '0.08221'
$ws.Range('E10').Value = '  +3.88%  '
$ws.Range('D11').Value = "'" + '1.037'
$ws.Range('E11').Value = '  +7.01%  '
$ws.Range('D12').Value = "'" + '22.76'
$ws.Range('E12').Value = '  +6.58%  '
$ws.Range('D13').Value = '1.958.82'
$ws.Range('E13').Value = '  +5.41%  '
$ws.Range('D14').Value = "'" + '6.152'
$ws.Range('E14').Value = '  +4.25%  '
$ws.Range('D15').Value = "'" + '7.382'
$ws.Range('E15').Value = '  +3.08%  '
$ws.Range('D16').Value = "'" + '91.87'
$ws.Range('E16').Value = '  +1.75%  '
$ws.Range('E17').Value = '  -0.62%  '
$ws.Range('D18').Value = "'" + '0.00001058'
$ws.Range('E18').Value = '  +2.85%  '
$ws.Range('E19').Value = '  +0.90%  '
$ws.Range('D20').Value = "'" + '18.04'
$ws.Range('E20').Value = '  +3.92%  '
$ws.Range('D21').Value = "'" + '1.000'
$ws.Range('E21').Value = '  -0.46%  '
$ws.Range('D22').Value = '29.750.50'
$ws.Range('E22').Value = '  +6.90%  '
$ws.Range('D23').Value = "'" + '5.582'
$ws.Range('E23').Value = '  +4.38%  '
$ws.Range('D24').Value = "'" + '11.25'
$ws.Range('E24').Value = '  +3.65%  '
$ws.Range('D25').Value = "'" + '2.291'
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('D26').Value = '2.186.81'
$ws.Range('E26').Value = '  +5.22%  '
$ws.Range('D27').Value = "'" + '161.26'
$ws.Range('E27').Value = '  +1.73%  '
$ws.Range('D28').Value = "'" + '20.21'
$ws.Range('E28').Value = '  +3.76%  '
$ws.Range('D29').Value = "'" + '2.172'
$ws.Range('E29').Value = '  +5.13%  '
$ws.Range('D30').Value = "'" + '5.658'
$ws.Range('E30').Value = '  +5.33%  '
$ws.Range('D31').Value = "'" + '122.88'
$ws.Range('E31').Value = '  +3.43%  '
$ws.Range('E32').Value = '  +6.35%  '
$ws.Range('D33').Value = "'" + '0.09654'
$ws.Range('E33').Value = '  +2.46%  '
$ws.Range('E34').Value = '  +10.97%  '
$ws.Range('D35').Value = "'" + '3.684'
$ws.Range('E35').Value = '  +2.61%  '
$ws.Range('D37').Value = "'" + '0.06260'
$ws.Range('E37').Value = '  +3.89%  '
$ws.Range('D38').Value = "'" + '0.02317'
$ws.Range('E38').Value = '  +4.86%  '
$ws.Range('D39').Value = "'" + '8.529'
$ws.Range('E39').Value = '  +3.17%  '
$ws.Range('E40').Value = '  +2.37%  '
$ws.Range('D41').Value = "'" + '0.6080'
$ws.Range('E41').Value = '  +4.56%  '
$ws.Range('D42').Value = "'" + '10.72'
$ws.Range('E42').Value = '  +6.42%  '
$ws.Range('B43').Value = 'Frax'
$ws.Range('C43').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D43').Value = "'" + '1.001'
$ws.Range('E43').Value = '  -0.43%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').Value = "'" + '0.1896'
$ws.Range('E44').Value = '  +2.76%  '
$ws.Range('E45').Value = '  +32.27%  '
$ws.Range('D46').Value = "'" + '1.274'
$ws.Range('E46').Value = '  -0.50%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'" + '12.52'
$ws.Range('E47').Value = '  +4.89%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').Value = "'" + '0.5704'
$ws.Range('E48').Value = '  +4.56%  '
$ws.Range('D49').Value = "'" + '0.07427'
$ws.Range('E49').Value = '  +8.44%  '
$ws.Range('D50').Value = "'" + '1.990'
$ws.Range('E50').Value = '  +3.04%  '
$ws.Range('D51').Value = "'" + '113.05'
$ws.Range('E51').Value = '  +2.05%  '
